$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update responsibility_options for the "Geometry-building-spins" and
# "Geometry-building-roomPlan" tasks (column E, rows 7 and 8).
$ws.Range("E7").Value = "IM-operators, auto"
$ws.Range("E8").Value = "IM-operators, client"

# Match the author's final selection/scroll position (topLeftCell B1,
# active cell E8).
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E8").Select()
